# Updated cryptos list: refresh the Price (column D) and Volume(1h)
# change (column E) figures for each coin row on the active sheet.
#
# Column D sometimes holds values that look numeric ("55.79", "240.52", ...)
# but the source data models them as plain text (e.g. "36.101.24" uses dots
# as thousands separators and would not round-trip as a number). Writing the
# string straight into `.Value` lets Excel's smart-parsing coerce it into a
# real number, so we briefly force a Text number format, assign the literal
# string, then clear the format again so the cell is left exactly as it was
# (General/no explicit style) but still holds a text value.
# Column E values are already padded with spaces around the percent sign
# ("  -2.03%  "), so Excel never treats them as numbers and they can be
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '36.031.29'; E = '  -2.03%  ' },
    @{ Row = 3; D = '1.931.93'; E = '  -4.84%  ' },
    @{ Row = 4; D = $null; E = '  +0.16%  ' },
    @{ Row = 5; D = '240.52'; E = '  -3.64%  ' },
    @{ Row = 6; D = $null; E = '  -5.88%  ' },
    @{ Row = 7; D = $null; E = '  +0.02%  ' },
    @{ Row = 8; D = '55.79'; E = '  -11.79%  ' },
    @{ Row = 9; D = '0.362'; E = '  -7.92%  ' },
    @{ Row = 10; D = '55.21'; E = '  -5.09%  ' },
    @{ Row = 11; D = '0.0816'; E = '  +2.25%  ' },
    @{ Row = 12; D = $null; E = '  -0.71%  ' },
    @{ Row = 13; D = '0.817'; E = '  -8.64%  ' },
    @{ Row = 14; D = '2.215.03'; E = '  -4.81%  ' },
    @{ Row = 15; D = '20.81'; E = '  -11.78%  ' },
    @{ Row = 16; D = '13.20'; E = '  -8.27%  ' },
    @{ Row = 17; D = '5.18'; E = '  -7.12%  ' },
    @{ Row = 18; D = '1.926.93'; E = '  -4.90%  ' },
    @{ Row = 19; D = '35.951.91'; E = '  -2.20%  ' },
    @{ Row = 20; D = '69.24'; E = '  -4.59%  ' },
    @{ Row = 21; D = $null; E = '  -3.73%  ' },
    @{ Row = 22; D = '227.36'; E = '  -4.12%  ' },
    @{ Row = 23; D = $null; E = '  -8.45%  ' },
    @{ Row = 24; D = $null; E = '  +0.02%  ' },
    @{ Row = 25; D = '2.43'; E = '  -4.17%  ' },
    @{ Row = 26; D = $null; E = '  -2.45%  ' },
    @{ Row = 27; D = '9.26'; E = '  -7.57%  ' },
    @{ Row = 28; D = '162.53'; E = '  +1.22%  ' },
    @{ Row = 29; D = $null; E = '  -6.00%  ' },
    @{ Row = 30; D = $null; E = '  -15.72%  ' },
    @{ Row = 31; D = '0.116'; E = '  -3.46%  ' },
    @{ Row = 32; D = '1.12'; E = '  -5.42%  ' },
    @{ Row = 33; D = '4.64'; E = '  -8.86%  ' },
    @{ Row = 34; D = '0.0619'; E = '  -0.99%  ' },
    @{ Row = 35; D = '4.24'; E = '  -6.00%  ' },
    @{ Row = 36; D = $null; E = '  +0.22%  ' },
    @{ Row = 37; D = '5.99'; E = '  -8.65%  ' },
    @{ Row = 38; D = '1.80'; E = '  -1.84%  ' },
    @{ Row = 39; D = '2.13'; E = '  -11.34%  ' },
    @{ Row = 40; D = '2.83'; E = '  -14.07%  ' },
    @{ Row = 41; D = '0.0960'; E = '  -4.77%  ' },
    @{ Row = 42; D = '2.88'; E = '  -2.04%  ' },
    @{ Row = 43; D = '1.16'; E = '  -7.76%  ' },
    @{ Row = 44; D = '0.0206'; E = '  -4.42%  ' },
    @{ Row = 45; D = $null; E = '  -9.41%  ' },
    @{ Row = 46; D = $null; E = '  -10.14%  ' },
    @{ Row = 47; D = '7.26'; E = '  -5.80%  ' },
    @{ Row = 48; D = '1.328.63'; E = '  -2.81%  ' },
    @{ Row = 49; D = '86.88'; E = '  -7.73%  ' },
    @{ Row = 50; D = '2.80'; E = '  -3.41%  ' },
    @{ Row = 51; D = '45.47'; E = '  +0.01%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.ClearFormats()
    }

    $ws.Range("E$row").Value = $u.E
}
